# EPS v3.3.1 -> v3.4.2 update for CApULAbIFM CO2 Abated per Unit Land Area by Impr For Mgmt.xlsx
#
# Summary of changes:
#  - Remove the "Calculations" worksheet (no longer used; its sole output
#    is replaced by a literal-derived formula on the CApULAbIFM sheet).
#  - On "CApULAbIFM": B2's formula now computes 1.5*10^6 directly instead
#    of referencing Calculations!A6.
#  - On "About": update the source citation from the old EPA report
#    citation (U.S. EPA / GHG report / URL / page ref / year) to a single
#    short note "consultation with American Forest Foundation"; remove the
#    now-unused hyperlink and the now-empty trailing rows.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. CApULAbIFM sheet: replace the Calculations!A6 reference with a
#    direct literal formula, since the Calculations sheet is going away.
# ---------------------------------------------------------------------
$calcSheet = $wb.Worksheets.Item("CApULAbIFM")
$calcSheet.Range("B2").Formula = "=1.5*10^6"
$calcSheet.Activate()
$calcSheet.Range("B3").Select()

# ---------------------------------------------------------------------
# 2. About sheet: rewrite the source citation block.
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop the hyperlink that pointed at the old EPA PDF (cell B6).
foreach ($h in @($about.Hyperlinks)) {
    $h.Delete()
}

# New, shorter citation text replacing "U.S. EPA".
$about.Range("B3").Value = "consultation with American Forest Foundation"

# Remove the now-unneeded trailing detail rows (old rows 6, 7, 9), then
# collapse the gap left behind so the final blank spacer row lands at
# row 6 (it keeps its original, non-custom 14.45pt row height).
$about.Rows("9").Delete()
$about.Rows("7").Delete()
$about.Rows("6").Delete()
$about.Rows("6").Delete()

# Old row 4 (the "2005" date row) is no longer needed at all.
$about.Range("A4:B4").Clear()

# Old row 5's label text is gone too; only the formatted, empty A5 cell
# remains (its bold font is cleared).
$about.Range("B5").Clear()
$about.Range("A5").Font.Bold = $false

$about.Activate()
$about.Range("C17").Select()

# ---------------------------------------------------------------------
# 3. Remove the now-obsolete "Calculations" worksheet entirely.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Calculations").Delete()

Write-Output "edit complete"
